# Edit Word document per commit:
# "TFS 13809 - New Submission: Changing Program/Behavior resets Direct/Indirect btn;"
#
# 1) Append a new bold bullet-style paragraph to the "Description" table
#    cell listing change items.
# 2) Append a new row to the Date/Change Description/Author changelog table.
# 3) Bump the Changeset number referenced in the Implementation Steps table.

$d = $word.ActiveDocument

# --- Edit 1: Description table (first table) -----------------------------
$descTable = $d.Tables.Item(1)
$descCell = $descTable.Cell(1, 2)
$r = $descCell.Range
$r.Find.Execute("label;", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter("`rTFS 13808 – New Submission: Changing Program resets Direct/Indirect btn;")

# --- Edit 2: Changelog table (second table) -------------------------------
$logTable = $d.Tables.Item(2)
$logTable.Rows.Add() | Out-Null
$newRowIndex = $logTable.Rows.Count
$logTable.Cell($newRowIndex, 1).Range.Text = "03/13/2019"
$logTable.Cell($newRowIndex, 2).Range.Text = "TSF 13808 – New Submission: Changing Program resets Direct/Indirect btn;"
$logTable.Cell($newRowIndex, 3).Range.Text = "Lili Huang"

# --- Edit 3: Changeset number bump ----------------------------------------
$d.Content.Find.Execute("41873", $true, $false, $false, $false, $false, $true, 1, $false, "41887", 2) | Out-Null

Write-Host "Edits applied."
